# Updates crypto price/volume cells to match the latest scrape.
# Numeric-looking "Price" strings (single-dot decimals) are written with a
# leading apostrophe so Excel stores them as literal text (matching the
# original inlineStr cells) instead of auto-converting to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.942.32"
$ws.Range("E2").Value = "  +4.32%  "
$ws.Range("D3").Value = "3.247.76"
$ws.Range("E3").Value = "  +2.48%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'577.80"
$ws.Range("E5").Value = "  +3.22%  "
$ws.Range("D6").Value = "'176.57"
$ws.Range("E6").Value = "  +2.73%  "
$ws.Range("D7").Value = "'0.605"
$ws.Range("E7").Value = "  +0.56%  "
$ws.Range("D9").Value = "3.247.47"
$ws.Range("E9").Value = "  +2.46%  "
$ws.Range("E10").Value = "  +4.36%  "
$ws.Range("E11").Value = "  +1.48%  "
$ws.Range("D12").Value = "'0.407"
$ws.Range("E12").Value = "  +3.02%  "
$ws.Range("D13").Value = "3.814.02"
$ws.Range("E13").Value = "  +2.62%  "
$ws.Range("E14").Value = "  +1.58%  "
$ws.Range("D15").Value = "'27.81"
$ws.Range("E15").Value = "  +0.97%  "
$ws.Range("D16").Value = "66.927.35"
$ws.Range("E16").Value = "  +4.41%  "
$ws.Range("E17").Value = "  +3.41%  "
$ws.Range("D18").Value = "3.248.09"
$ws.Range("E18").Value = "  +2.53%  "
$ws.Range("D19").Value = "'5.81"
$ws.Range("E19").Value = "  +2.51%  "
$ws.Range("D20").Value = "'13.34"
$ws.Range("E20").Value = "  +2.47%  "
$ws.Range("D21").Value = "'368.56"
$ws.Range("E21").Value = "  +4.53%  "
$ws.Range("E22").Value = "  +4.34%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "'70.46"
$ws.Range("E24").Value = "  +1.92%  "
$ws.Range("D25").Value = "'0.506"
$ws.Range("E25").Value = "  +1.02%  "
$ws.Range("D26").Value = "3.384.20"
$ws.Range("E26").Value = "  +2.47%  "
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("D28").Value = "'9.78"
$ws.Range("E28").Value = "  +3.55%  "
$ws.Range("E29").Value = "  +1.64%  "
$ws.Range("E30").Value = "  +0.14%  "
$ws.Range("E31").Value = "  +5.05%  "
$ws.Range("D32").Value = "'5.64"
$ws.Range("E32").Value = "  +0.78%  "
$ws.Range("D33").Value = "'22.48"
$ws.Range("E33").Value = "  +1.81%  "
$ws.Range("D35").Value = "'172.05"
$ws.Range("E35").Value = "  +9.59%  "
$ws.Range("E36").Value = "  +2.49%  "
$ws.Range("D37").Value = "'6.76"
$ws.Range("E37").Value = "  +2.26%  "
$ws.Range("E38").Value = "  +4.72%  "
$ws.Range("D39").Value = "'0.851"
$ws.Range("E39").Value = "  +6.97%  "
$ws.Range("E40").Value = "  +9.91%  "
$ws.Range("D41").Value = "'26.84"
$ws.Range("E41").Value = "  +3.09%  "
$ws.Range("E42").Value = "  +2.07%  "
$ws.Range("D43").Value = "'6.42"
$ws.Range("E43").Value = "  +6.22%  "
$ws.Range("D44").Value = "2.717.99"
$ws.Range("E44").Value = "  +2.54%  "
$ws.Range("D45").Value = "'4.29"
$ws.Range("E45").Value = "  +3.29%  "
$ws.Range("D46").Value = "'40.43"
$ws.Range("E46").Value = "  +4.23%  "
$ws.Range("D47").Value = "'339.15"
$ws.Range("E47").Value = "  +3.95%  "
$ws.Range("E48").Value = "  +3.31%  "
$ws.Range("D49").Value = "'24.56"
$ws.Range("E49").Value = "  +4.18%  "
$ws.Range("E50").Value = "  +3.26%  "
$ws.Range("E51").Value = "  +2.58%  "
